# "ER atualizado e tabelas mapeadas na pasta 'tabelas novas'"
#
# The old mapping row "cidade" -> "city_id" (row 8) is no longer needed,
# so remove that whole row from the mapping table. Everything below it
# (including the "ativo" -> "actived" row) shifts up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row that holds the "cidade" / "city_id" mapping.
$ws.Rows(8).Delete()

# Leave the selection where the author left off editing.
$ws.Range("B17").Select()
